# Auto-generated Excel COM-interop script
# Applies numeric corrections to columns H-N across several sheets
# in the Ragnarok_Profits workbook (market-profit recalculation refresh).

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2001579.8
$ws.Range("I2").Value = 3000
$ws.Range("J2").Value = 2501224.8
$ws.Range("K2").Value = 3000
$ws.Range("L2").Value = 2501224.8
$ws.Range("M2").Value = -2887
$ws.Range("N2").Value = -2501450.8
$ws.Range("H33").Value = 322.5
$ws.Range("I33").Value = 296.77777
$ws.Range("K33").Value = 296.77777
$ws.Range("M33").Value = -67.77776999999998
$ws.Range("H103").Value = 62501460
$ws.Range("I103").Value = 922.5
$ws.Range("J103").Value = 125002000
$ws.Range("K103").Value = 2767.5
$ws.Range("L103").Value = 375006000
$ws.Range("M103").Value = -2181.5
$ws.Range("N103").Value = -375007172
$ws.Range("H112").Value = 3046.476
$ws.Range("I112").Value = 1295.5
$ws.Range("J112").Value = 3458.4707
$ws.Range("K112").Value = 3886.5
$ws.Range("L112").Value = 10375.4121
$ws.Range("M112").Value = -2778.5
$ws.Range("N112").Value = -12591.4121
$ws.Range("H137").Value = 1873.5128
$ws.Range("I137").Value = 1754.9032
$ws.Range("K137").Value = 5264.7096
$ws.Range("M137").Value = -2714.7096
$ws.Range("H138").Value = 7296.3706
$ws.Range("J138").Value = 6916.52
$ws.Range("L138").Value = 20749.56
$ws.Range("N138").Value = -31029.56

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4831.2905
$ws.Range("I2").Value = 5223.8423
$ws.Range("J2").Value = 4209.75
$ws.Range("K2").Value = 5223.8423
$ws.Range("L2").Value = 4209.75
$ws.Range("M2").Value = -5110.8423
$ws.Range("N2").Value = -4435.75
$ws.Range("H32").Value = 5337.4707
$ws.Range("I32").Value = 5337.4707
$ws.Range("K32").Value = 5337.4707
$ws.Range("M32").Value = -5050.4707
$ws.Range("H45").Value = 3843.7
$ws.Range("I45").Value = 3505.75
$ws.Range("K45").Value = 3505.75
$ws.Range("M45").Value = -3128.75
$ws.Range("H74").Value = 3092.5386
$ws.Range("I74").Value = 1820.5
$ws.Range("K74").Value = 1820.5
$ws.Range("M74").Value = -946.5
$ws.Range("H77").Value = 3092.5386
$ws.Range("I77").Value = 1820.5
$ws.Range("K77").Value = 9102.5
$ws.Range("M77").Value = -4734.5
$ws.Range("H116").Value = 4831.2905
$ws.Range("I116").Value = 5223.8423
$ws.Range("J116").Value = 4209.75
$ws.Range("K116").Value = 5223.8423
$ws.Range("L116").Value = 4209.75
$ws.Range("M116").Value = -2929.8423
$ws.Range("N116").Value = -8797.75
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").Value = $null
$ws.Range("H125").Value = 71111
$ws.Range("J125").Value = 71111
$ws.Range("L125").Value = 71111
$ws.Range("N125").Value = -80951
$ws.Range("H132").Value = 2781809.2
$ws.Range("I132").Value = 3958.5757
$ws.Range("K132").Value = 11875.7271
$ws.Range("M132").Value = -9345.7271

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4831.2905
$ws.Range("I3").Value = 5223.8423
$ws.Range("J3").Value = 4209.75
$ws.Range("K3").Value = 5223.8423
$ws.Range("L3").Value = 4209.75
$ws.Range("M3").Value = -5109.8423
$ws.Range("N3").Value = -4437.75
$ws.Range("H22").Value = 2562.5
$ws.Range("I22").Value = 1625.5
$ws.Range("K22").Value = 1625.5
$ws.Range("M22").Value = -1452.5
$ws.Range("H94").Value = 3366.4
$ws.Range("I94").Value = 3392.6956
$ws.Range("K94").Value = 3392.6956
$ws.Range("M94").Value = -2941.6956
$ws.Range("H99").Value = 2587.8
$ws.Range("I99").Value = 2628.4443
$ws.Range("K99").Value = 2628.4443
$ws.Range("M99").Value = -1130.4443
$ws.Range("H105").Value = 808625.1
$ws.Range("I105").Value = 1432112.5
$ws.Range("K105").Value = 1432112.5
$ws.Range("M105").Value = -1430365.5
$ws.Range("H110").Value = 89993
$ws.Range("J110").Value = 89993
$ws.Range("L110").Value = 89993
$ws.Range("N110").Value = -98173
$ws.Range("H134").Value = 4547959
$ws.Range("I134").Value = 2595
$ws.Range("J134").Value = 33335264
$ws.Range("K134").Value = 7785
$ws.Range("L134").Value = 100005792
$ws.Range("M134").Value = -5250
$ws.Range("N134").Value = -100010862

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27780954
$ws.Range("I31").Value = 38464590
$ws.Range("K31").Value = 38464590
$ws.Range("M31").Value = -38464295
$ws.Range("H34").Value = 27780954
$ws.Range("I34").Value = 38464590
$ws.Range("K34").Value = 38464590
$ws.Range("M34").Value = -38464388
$ws.Range("H62").Value = 4249.5
$ws.Range("J62").Value = 4499.5
$ws.Range("L62").Value = 4499.5
$ws.Range("N62").Value = -5747.5
$ws.Range("H65").Value = 4249.5
$ws.Range("J65").Value = 4499.5
$ws.Range("L65").Value = 22497.5
$ws.Range("N65").Value = -28737.5
$ws.Range("H132").Value = 2721.7407
$ws.Range("I132").Value = 2624.5
$ws.Range("J132").Value = 3499.6667
$ws.Range("K132").Value = 7873.5
$ws.Range("L132").Value = 10499.0001
$ws.Range("M132").Value = -5343.5
$ws.Range("N132").Value = -15559.0001
$ws.Range("H134").Value = 3247.8333
$ws.Range("I134").Value = 3179.4546
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 9538.363799999999
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -7003.363799999999
$ws.Range("N134").Value = -17070

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1387.15
$ws.Range("I107").Value = 1385.0588
$ws.Range("J107").Value = 1399
$ws.Range("K107").Value = 1385.0588
$ws.Range("L107").Value = 1399
$ws.Range("M107").Value = 534.9412
$ws.Range("N107").Value = -5239
$ws.Range("H126").Value = 1462.5
$ws.Range("I126").Value = 1612
$ws.Range("J126").Value = 1014
$ws.Range("K126").Value = 4836
$ws.Range("L126").Value = 3042
$ws.Range("M126").Value = -2366
$ws.Range("N126").Value = -7982
$ws.Range("H132").Value = 16671243
$ws.Range("I132").Value = 5492
$ws.Range("K132").Value = 16476
$ws.Range("M132").Value = -13946

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3247.5
$ws.Range("I46").Value = 2495
$ws.Range("J46").Value = 4000
$ws.Range("K46").Value = 2495
$ws.Range("L46").Value = 4000
$ws.Range("M46").Value = -2307
$ws.Range("N46").Value = -4376
$ws.Range("H68").Value = 1986043.8
$ws.Range("I68").Value = 3206301.5
$ws.Range("J68").Value = 3125
$ws.Range("K68").Value = 3206301.5
$ws.Range("L68").Value = 3125
$ws.Range("M68").Value = -3205552.5
$ws.Range("N68").Value = -4623
$ws.Range("H71").Value = 1986043.8
$ws.Range("I71").Value = 3206301.5
$ws.Range("J71").Value = 3125
$ws.Range("K71").Value = 16031507.5
$ws.Range("L71").Value = 15625
$ws.Range("M71").Value = -16027763.5
$ws.Range("N71").Value = -23113
$ws.Range("H122").Value = 3619.244
$ws.Range("I122").Value = 3310.639
$ws.Range("K122").Value = 9931.917000000001
$ws.Range("M122").Value = -7481.917000000001
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").Value = $null
$ws.Range("H136").Value = 4438.5
$ws.Range("I136").Value = 2834.6667
$ws.Range("J136").Value = 9250
$ws.Range("K136").Value = 8504.000100000001
$ws.Range("L136").Value = 27750
$ws.Range("M136").Value = -5954.000100000001
$ws.Range("N136").Value = -32850

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 27958.334
$ws.Range("I15").Value = 24450
$ws.Range("J15").Value = 34975
$ws.Range("K15").Value = 24450
$ws.Range("L15").Value = 34975
$ws.Range("M15").Value = -24162
$ws.Range("N15").Value = -35551
$ws.Range("H126").Value = 6368.6665
$ws.Range("I126").Value = 7502.9
$ws.Range("K126").Value = 22508.7
$ws.Range("M126").Value = -20038.7
$ws.Range("H132").Value = 1001903.9
$ws.Range("I132").Value = 1754.875
$ws.Range("K132").Value = 5264.625
$ws.Range("M132").Value = -2734.625
$ws.Range("H136").Value = 632474.25
$ws.Range("I136").Value = 8892.23
$ws.Range("K136").Value = 26676.69
$ws.Range("M136").Value = -24126.69

Write-Output "Applied 210 cell updates and 2 cell clears across 7 sheets."